$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.257.05"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "2.485.22"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.08"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.05"
$ws.Range("E6").Value = "  +2.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.521"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.539"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.14"
$ws.Range("E10").Value = "  +3.08%  "
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.49"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").Value = "2.873.39"
$ws.Range("D16").Value = "2.486.00"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "47.180.21"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.42"
$ws.Range("E19").Value = "  +5.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.63"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.79"
$ws.Range("E21").Value = "  +17.14%  "
$ws.Range("D22").Value = "0.0₃0941"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "246.19"
$ws.Range("E24").Value = "  -2.07%  "
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.66"
$ws.Range("E27").Value = "  -2.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.30"
$ws.Range("E28").Value = "  +4.30%  "
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.139"
$ws.Range("E30").Value = "  +3.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.55"
$ws.Range("E31").Value = "  -1.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.84"
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.29"
$ws.Range("E33").Value = "  +1.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.31"
$ws.Range("E34").Value = "  -1.30%  "
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.76"
$ws.Range("E37").Value = "  +3.07%  "
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.94"
$ws.Range("E40").Value = "  +8.88%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.21"
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.03"
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0296"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "1.996.62"
$ws.Range("E45").Value = "  +1.62%  "
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.01"
$ws.Range("E47").Value = "  -4.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.79"
$ws.Range("E48").Value = "  -2.43%  "
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.11"
$ws.Range("E50").Value = "  -4.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.73"
$ws.Range("E51").Value = "  +3.33%  "
